$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 02:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1620204
$ws.Range("C4").Value = 27481
$ws.Range("D4").Value = 382183
$ws.Range("E4").Value = 1141739
$ws.Range("G4").Value = 1346
$ws.Range("H4").Value = 96282

# Canada (row 17)
$ws.Range("B17").Value = 81324
$ws.Range("C17").Value = 1182
$ws.Range("D17").Value = 41715
$ws.Range("E17").Value = 33457
$ws.Range("G17").Value = 121
$ws.Range("H17").Value = 6152

# Panama (row 50)
$ws.Range("B50").Value = 10116
$ws.Range("C50").Value = 139
$ws.Range("D50").Value = 6245
$ws.Range("E50").Value = 3580
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 291

# Cabo Verde (row 140)
$ws.Range("D140").Value = 95
$ws.Range("E140").Value = 258

# Bahamas (row 170)
$ws.Range("D170").Value = 44
$ws.Range("E170").Value = 42
